$d = $word.ActiveDocument

# The "Non-functional requirements" list (numId 37, ilvl 2) has five
# heading/body pairs: Availability, Correctness, Maintainability,
# Reusability, Portability. Each heading paragraph currently inherits its
# 2160-twip (108pt) left indent from the numbering definition, and each
# body paragraph has an explicit 2160-twip (108pt) left indent. Both need
# to become 1800 twips (90pt): an explicit override added to the heading
# paragraphs, and the explicit value changed on the body paragraphs.

$targets = @("Availability", "Correctness", "Maintainability", "Reusability", "Portability")

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.Trim()

    if ($targets -contains $text) {
        # Heading paragraph (e.g. "Availability") -- set its left indent,
        # which adds the explicit <w:ind w:left="1800"/> override.
        $p.Format.LeftIndent = 90

        # The very next paragraph is the body text with the 2160 indent.
        $next = $p.Next()
        if ($next -ne $null -and $next.Format.LeftIndent -eq 108) {
            $next.Format.LeftIndent = 90
        }
    }
}
